$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.946.73"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "2.354.14"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.685"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.75"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.64%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.613"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +13.85%  "
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.25"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "33.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +16.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +10.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.107"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("D15").Value = "2.702.38"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("E16").Value = "  -2.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.912"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.99%  "
$ws.Range("D18").Value = "2.319.54"
$ws.Range("E18").Value = "  -2.31%  "
$ws.Range("D19").Value = "43.865.91"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "258.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.64%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +15.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("E32").Value = "  -4.67%  "
$ws.Range("E33").Value = "  +4.17%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.41%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0754"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0281"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.215"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +18.89%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.108"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.12%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.17%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.47%  "
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.14%  "
$ws.Range("E48").Value = "  +2.48%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "100.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.45%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.43%  "
